$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0
$ws.Range("B2").Value = -0.107670151431414
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = 0.2053876115891933
$ws.Range("E2").Value = -0.000870470479885202
$ws.Range("G2").Value = 0
$ws.Range("J2").Value = -0
$ws.Range("K2").Value = -0.01825145604069086
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.2064948814758124
$ws.Range("N2").Value = 0.00208025141448778
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = -0.09778338212392844
$ws.Range("V2").Value = 0.01432430740601758
$ws.Range("W2").Value = -0.04174952113749071
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = -0.07844340024108128
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = -0.01123724989979894
$ws.Range("AF2").Value = 0.00467995902030148
$ws.Range("AG2").Value = -0
$ws.Range("AI2").Value = -0
$ws.Range("AJ2").Value = 0
$ws.Range("AK2").Value = -0
$ws.Range("AL2").Value = -0.03477317437924476
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.02975591359928845
$ws.Range("AO2").Value = 0.06894031625633901
$ws.Range("AR2").Value = -0
$ws.Range("AT2").Value = 0
$ws.Range("AU2").Value = -0.1504158949591624
$ws.Range("AW2").Value = 0.06762809990804736
$ws.Range("AX2").Value = 0.01316411438484101
$ws.Range("AY2").Value = -0
$ws.Range("BC2").Value = -0
$ws.Range("BD2").Value = -0.03736677967854636
$ws.Range("BF2").Value = 0.07026386262018187
$ws.Range("BG2").Value = 0.02677008819260617
$ws.Range("BI2").Value = -0
$ws.Range("BJ2").Value = -0
$ws.Range("BL2").Value = 0
$ws.Range("BM2").Value = 0.01906614199068762
$ws.Range("BO2").Value = -0.05327664047297926
$ws.Range("BP2").Value = -0.06915173979540963
$ws.Range("BU2").Value = 0
$ws.Range("BV2").Value = -0.01804737642947697
$ws.Range("BX2").Value = 0.01512898727912453
$ws.Range("BY2").Value = -0.005352203034577101
$ws.Range("BZ2").Value = -0
$ws.Range("CB2").Value = 0
$ws.Range("CD2").Value = 0
$ws.Range("CE2").Value = 0.0336355770062513
$ws.Range("CG2").Value = -0.03569317361011079
$ws.Range("CH2").Value = 0.0251507005628732
$ws.Range("CJ2").Value = -0
$ws.Range("CM2").Value = -0
$ws.Range("CN2").Value = -0.01586600892987506
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = 0.02594067695371425
$ws.Range("CQ2").Value = 0.05198710927025178
$ws.Range("CT2").Value = 0
$ws.Range("CU2").Value = -0
$ws.Range("CV2").Value = -0
$ws.Range("CW2").Value = 0.04516854268119422
$ws.Range("CY2").Value = -0.03815737946181457
$ws.Range("CZ2").Value = 0.01954564062428951
$ws.Range("DB2").Value = 0
$ws.Range("DE2").Value = -0
$ws.Range("DF2").Value = 0.04358123142453846
$ws.Range("DH2").Value = 0.01717012150273894
$ws.Range("DI2").Value = 0.04988333764285281
$ws.Range("DJ2").Value = 0
$ws.Range("DK2").Value = -0
$ws.Range("DL2").Value = -0
$ws.Range("DN2").Value = 0
$ws.Range("DO2").Value = -0.02121619411584265
$ws.Range("DQ2").Value = 0.03359901195050785
$ws.Range("DR2").Value = -0.01385797544263802
$ws.Range("DS2").Value = -0
$ws.Range("DW2").Value = 0
$ws.Range("DX2").Value = -0.05212903656296591
$ws.Range("DY2").Value = -0
$ws.Range("DZ2").Value = -0.004971154324998677
$ws.Range("EA2").Value = -0.02718932500684638
$ws.Range("EB2").Value = 0
$ws.Range("EF2").Value = -0
$ws.Range("EG2").Value = 0.05015491900709527
$ws.Range("EI2").Value = 0.07613886383889099
$ws.Range("EJ2").Value = -0.0199689068547195
$ws.Range("EO2").Value = 0
$ws.Range("EP2").Value = 0.04747992163333337
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = -0.03421357883596667
$ws.Range("ES2").Value = 0.04056784176355664
$ws.Range("ET2").Value = 0
$ws.Range("EU2").Value = -0
$ws.Range("EV2").Value = 0
$ws.Range("EX2").Value = 0
$ws.Range("EY2").Value = 0.047632734695192
$ws.Range("FA2").Value = -0.03127982317017986
$ws.Range("FB2").Value = 0.02441504420193668
$ws.Range("FD2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = -0
$ws.Range("FH2").Value = -0.0002063141115742818
$ws.Range("FI2").Value = 0
$ws.Range("FJ2").Value = -0.01377775761390556
$ws.Range("FK2").Value = -0.001886070509032845
$ws.Range("FL2").Value = -0
$ws.Range("FP2").Value = -0
$ws.Range("FQ2").Value = -0.01316388935840617
$ws.Range("FR2").Value = -0
$ws.Range("FS2").Value = -0.0009316850821103231
$ws.Range("FT2").Value = 0.0192548302470746
$ws.Range("FV2").Value = -0
$ws.Range("FW2").Value = -0
$ws.Range("FY2").Value = 0
$ws.Range("FZ2").Value = -0.03444598377532641
$ws.Range("GA2").Value = -0
$ws.Range("GB2").Value = 0.04121121055052257
$ws.Range("GD2").Value = 0
